$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 387 (old rows 387-433 shift down to 388-434).
$ws.Rows("387:387").Insert()

# Populate the newly inserted row 387 with the new record's data.
# Values mirror the old row 387 (Segunda / Ecuador / 14-unit box) except for
# the Fecha (D) and Volumen (M), which are the new data points.
$ws.Cells.Item(387, 1).Value2 = 5
$ws.Cells.Item(387, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(387, 3).Value2 = "Maule"
$ws.Cells.Item(387, 4).Value2 = 45142
$ws.Cells.Item(387, 5).Value2 = 7
$ws.Cells.Item(387, 6).Value2 = "Fruta"
$ws.Cells.Item(387, 7).Value2 = 100108
$ws.Cells.Item(387, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(387, 9).Value2 = 100108005
$ws.Cells.Item(387, 10).Value2 = "Piña"
$ws.Cells.Item(387, 11).Value2 = "Caramelo"
$ws.Cells.Item(387, 12).Value2 = "Segunda"
$ws.Cells.Item(387, 13).Value2 = 200
$ws.Cells.Item(387, 14).Value2 = 20000
$ws.Cells.Item(387, 15).Value2 = 20000
$ws.Cells.Item(387, 16).Value2 = 20000
$ws.Cells.Item(387, 17).Value2 = "`$/caja 14 unidades"
$ws.Cells.Item(387, 18).Value2 = "Ecuador"
$ws.Cells.Item(387, 19).Value2 = 1429
$ws.Cells.Item(387, 20).Value2 = 14
